# update(CauHoi): import CauHoi (TN+Dien)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename sheets: Sheet5 -> DienTu, Sheet6 -> NoiTu
# ---------------------------------------------------------------------
$wsDienTu = $wb.Worksheets.Item(3)
$wsDienTu.Name = "DienTu"
$wsNoiTu = $wb.Worksheets.Item(4)
$wsNoiTu.Name = "NoiTu"

$wsCauHoi = $wb.Worksheets.Item(1)
$wsCauTraLoi = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# 2. CauHoi sheet (sheet1): insert a new column A ("MaCauHoi"), shift
#    existing columns right, and append a new "fill the blank" question.
# ---------------------------------------------------------------------

# Preserve exact widths of the (soon to be shifted) columns before the insert.
$wOld1 = $wsCauHoi.Columns.Item(1).Width
$wOld2 = $wsCauHoi.Columns.Item(2).Width
$wOld3 = $wsCauHoi.Columns.Item(3).Width
$wOld4 = $wsCauHoi.Columns.Item(4).Width

$wsCauHoi.Columns.Item(1).Insert()

$wsCauHoi.Columns.Item(2).Width = $wOld1
$wsCauHoi.Columns.Item(3).Width = $wOld2
$wsCauHoi.Columns.Item(4).Width = $wOld3
$wsCauHoi.Columns.Item(5).Width = $wOld4

# Header row
$wsCauHoi.Range("A1").Value = "MaCauHoi"

# New MaCauHoi numbering column (STT) for existing rows
$wsCauHoi.Range("A2").Value = 1
$wsCauHoi.Range("A3").Value = 2
$wsCauHoi.Range("A4").Value = 3

$wsCauHoi.Range("A2:A4").HorizontalAlignment = -4108
$wsCauHoi.Range("A2:A4").VerticalAlignment = -4108

# D4 (DoKho) changed from 3 -> 3 already correct, ensure LoaiCauHoi col values are style 1
$wsCauHoi.Range("C2:E4").HorizontalAlignment = -4108
$wsCauHoi.Range("C2:E4").VerticalAlignment = -4108

# New row 5: a "fill in the blank" (Dien tu) question
$wsCauHoi.Range("A5").Value = 4
$wsCauHoi.Range("A5").HorizontalAlignment = -4108
$wsCauHoi.Range("A5").VerticalAlignment = -4108

$wsCauHoi.Range("B5").Value = "Mo hinh (1) la mot quy trinh phat trien phan mem, trong do cac giai doan nhu yeu cau, thiet ke, trien khai va kiem thu duoc thuc hien theo thu tu tuan tu. (Thac nuoc)"
$wsCauHoi.Range("B5").HorizontalAlignment = -4131
$wsCauHoi.Range("B5").VerticalAlignment = -4108

$wsCauHoi.Range("C5").Value = 3
$wsCauHoi.Range("D5").Value = 2
$wsCauHoi.Range("C5:D5").HorizontalAlignment = -4108
$wsCauHoi.Range("C5:D5").VerticalAlignment = -4108

$wsCauHoi.Range("E5").Value = "Dien tu"
$wsCauHoi.Range("E5").HorizontalAlignment = -4108
$wsCauHoi.Range("E5").VerticalAlignment = -4108

$wsCauHoi.Range("B8").Value = ""

# ---------------------------------------------------------------------
# 3. DienTu sheet (sheet3): populate with ViTri / DapAnText data
# ---------------------------------------------------------------------
$wsDienTu.Range("A1").Value = "MaCauHoi"
$wsDienTu.Range("B1").Value = "ViTri"
$wsDienTu.Range("C1").Value = "DapAnText"

$wsDienTu.Range("A2").Value = 4
$wsDienTu.Range("B2").Value = 1
$wsDienTu.Range("C2").Value = "Thac nuoc"

$wsDienTu.Range("A3").Value = 5
$wsDienTu.Range("B3").Value = 3
$wsDienTu.Range("C3").Value = "abc"

$wsDienTu.Columns.Item(1).ColumnWidth = 16.71
$wsDienTu.Columns.Item(2).ColumnWidth = 18.28
$wsDienTu.Columns.Item(3).ColumnWidth = 29.71

# ---------------------------------------------------------------------
# 4. Sheet selections / active sheet bookkeeping
# ---------------------------------------------------------------------
$wsCauTraLoi.Range("B25").Select()
$wsDienTu.Range("F6").Select()

$wsCauHoi.Activate()
$wsCauHoi.Range("D9").Select()
